# "borrado registro horas viejo" - fill in the weekly hours table on the
# first sheet (Hoja1) with the missing entries, and update the per-row /
# per-column totals (which recalc automatically from the formulas already
# in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Semana  1 (17/8) - row 3 : "Otros" block (L:O) was empty -----------
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.5
$ws.Range("N3").Value = 0.5
$ws.Range("O3").Value = 0.5

# --- Semana2 (24/8) - row 4 : "Verificacion/Testing" block (G:J) -------
$ws.Range("I4").Value = 1.5
$ws.Range("J4").Value = 1

# --- Semana 3 (31/8) - row 5 : "Verificacion/Testing" block (G:J) ------
$ws.Range("H5").Value = 0.5
$ws.Range("I5").Value = 1.5
$ws.Range("J5").Value = 1.5

# --- Semana 4 - row 6 : "Otros" block (L:O) -----------------------------
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 2

# --- row 9 : "Verificacion/Testing" (G:J) and "Otros" (L:O) blocks -----
$ws.Range("I9").Value = 4
$ws.Range("N9").Value = 1.5
$ws.Range("O9").Value = 1

# --- TOTAL column (V) : rows 3-5 now also add the "Investigacion" (K)
#     column into their total, so they get their own (non-shared) formula,
#     while rows 6-10 keep sharing the original "=SUM(U,P,F)" pattern.
$ws.Range("V3").Formula = "=SUM(U3,P3,F3,K3)"
$ws.Range("V4").Formula = "=SUM(U4,P4,F4,K4)"
$ws.Range("V5").Formula = "=SUM(U5,K5,P5,F5)"
$ws.Range("V6:V10").Formula = "=SUM(U6,P6,F6)"

# --- last active selection recorded on the sheet -----------------------
$ws.Range("V17").Select()
